# Apply the "output generated at 456a3b4" update to 广州-漫展信息.xlsx
# Sheets: 展览 (Exhibition), 演出 (Performance), 本地生活 (Local life), 全部类型 (All types)

$wb = $excel.ActiveWorkbook

$wsExpo  = $wb.Worksheets.Item("展览")
$wsShow  = $wb.Worksheets.Item("演出")
$wsLife  = $wb.Worksheets.Item("本地生活")
$wsAll   = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------
# Sheet "展览" (Exhibition) - update "想去人数" (want-to-go count) column F
# ---------------------------------------------------------------
$wsExpo.Range("F3").Value  = 970
$wsExpo.Range("F4").Value  = 99
$wsExpo.Range("F5").Value  = 0
$wsExpo.Range("F6").Value  = 480
$wsExpo.Range("F7").Value  = 78
$wsExpo.Range("F8").Value  = 1523
$wsExpo.Range("F9").Value  = 39230
$wsExpo.Range("G9").Value  = "已售罄"
$wsExpo.Range("F10").Value = 8302
$wsExpo.Range("F19").Value = 50
$wsExpo.Range("F22").Value = 213
$wsExpo.Range("F23").Value = 0
$wsExpo.Range("F27").Value = 550
$wsExpo.Range("F30").Value = 18
$wsExpo.Range("F36").Value = 173
$wsExpo.Range("F41").Value = 1011
$wsExpo.Range("F42").Value = 331
$wsExpo.Range("F43").Value = 72
$wsExpo.Range("F44").Value = 7

# ---------------------------------------------------------------
# Sheet "演出" (Performance) - update column F
# ---------------------------------------------------------------
$wsShow.Range("F2").Value  = 4
$wsShow.Range("F3").Value  = 190
$wsShow.Range("F5").Value  = 4363
$wsShow.Range("F7").Value  = 292
$wsShow.Range("F8").Value  = 0
$wsShow.Range("F9").Value  = 5
$wsShow.Range("F10").Value = 71
$wsShow.Range("F11").Value = 89
$wsShow.Range("F13").Value = 46
$wsShow.Range("F19").Value = 4361

# ---------------------------------------------------------------
# Sheet "本地生活" (Local life) - update column F
# ---------------------------------------------------------------
$wsLife.Range("F2").Value = 1750

# ---------------------------------------------------------------
# Sheet "全部类型" (All types) - update column F (and a block of rows
# that shifted because an older performance event dropped off the list
# and a new one was appended)
# ---------------------------------------------------------------
$wsAll.Range("F2").Value  = 1750
$wsAll.Range("F3").Value  = 424
$wsAll.Range("F6").Value  = 970
$wsAll.Range("F8").Value  = 914
$wsAll.Range("F9").Value  = 4
$wsAll.Range("F10").Value = 0

# Rows 11-15 on "全部类型" represent the same events as rows 2-8 on
# "演出"; the 07-19 "萤火虫" event disappeared from the feed and every
# subsequent row shifted up by one, with a brand-new "浪漫古典II" event
# appended at the end (row 15).

$wsAll.Range("B11").Value = "2024-07-20"
$wsAll.Range("C11").Value = "广州·冰兔2024线下live「过去和未来」"
$wsAll.Range("D11").Value = "恩宁路265号三层四层自编01 MAO Livehouse广州（永庆坊店）"
$wsAll.Range("E11").Value = "2024.07.20 20:00-07.20 22:00"
$wsAll.Range("F11").Value = 190
$wsAll.Range("G11").Value = 198
$wsAll.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=87546"
$wsAll.Range("I11").Value = "//i2.hdslb.com/bfs/openplatform/202406/2X09PE1a1718611339266.jpeg"

$wsAll.Range("B12").Value = "2024-07-20"
$wsAll.Range("C12").Value = "广州·跨越二次元ACG神级动漫世界巡回演唱会"
$wsAll.Range("D12").Value = "广州市荔湾区十甫路125号(上下九步行街内)2层 广州平安大戏院"
$wsAll.Range("E12").Value = "2024.07.20 19:30-07.20 21:10"
$wsAll.Range("F12").Value = 324
$wsAll.Range("G12").Value = 280
$wsAll.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=85353"
$wsAll.Range("I12").Value = "//i1.hdslb.com/bfs/openplatform/202405/4gACWbPh1715223804704.jpeg"

$wsAll.Range("B13").Value = "2024-07-21"
$wsAll.Range("C13").Value = "广州·昨日重现——唯美英文经典歌曲演唱会"
$wsAll.Range("D13").Value = "东风中路299号 广州中山纪念堂"
$wsAll.Range("E13").Value = "2024.07.21 19:30-07.21 21:30"
$wsAll.Range("F13").Value = 7
$wsAll.Range("G13").Value = 100
$wsAll.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=86802"
$wsAll.Range("I13").Value = "//i1.hdslb.com/bfs/openplatform/202405/DR8AvmXe1716802703006.jpeg"

$wsAll.Range("B14").Value = "2024-07-21"
$wsAll.Range("C14").Value = "广州·燃动!!高梨康治SUMMER LIVE TOUR IN CHINA 2024"
$wsAll.Range("D14").Value = "海珠同创汇东一街11号（上冲南约11-2） 声音共和Livehouse"
$wsAll.Range("E14").Value = "2024.07.21 14:30-07.21 16:00"
$wsAll.Range("F14").Value = 292
$wsAll.Range("G14").Value = 280
$wsAll.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=87034"
$wsAll.Range("I14").Value = "//i1.hdslb.com/bfs/openplatform/202406/LINsP2ui1717741701901.png"

$wsAll.Range("B15").Value = "2024-07-26"
$wsAll.Range("C15").Value = "广州·【早鸟8折】“浪漫古典Ⅱ”百年经典传世名曲烛光音乐会 "
$wsAll.Range("D15").Value = "广州市二沙岛晴波路33号  星海音乐厅（交响乐演奏厅）"
$wsAll.Range("E15").Value = "2024.07.26 20:00-07.26 21:30"
$wsAll.Range("F15").Value = 7
$wsAll.Range("G15").Value = 144
$wsAll.Range("H15").Value = "https://show.bilibili.com/platform/detail.html?id=87726"
$wsAll.Range("I15").Value = "//i0.hdslb.com/bfs/openplatform/202406/A8vhVlhn1717575084179.png"

$wsAll.Range("F16").Value = 8302
$wsAll.Range("F19").Value = 71
$wsAll.Range("F22").Value = 89
$wsAll.Range("F23").Value = 89
$wsAll.Range("F31").Value = 213
$wsAll.Range("F36").Value = 550
$wsAll.Range("F38").Value = 33
$wsAll.Range("F39").Value = 18
$wsAll.Range("F40").Value = 377
$wsAll.Range("F46").Value = 1011
$wsAll.Range("F49").Value = 12
